$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts old rows 8-24 down to 9-25,
# carrying their values/styles with them).
$ws.Rows.Item(8).Insert()

# The new row 8 should look like the "Season end" row (row 7) for the
# B/C/D formatting (blank trigger cell + date-formatted values), so copy
# that row's formatting across before filling in the new content.
$ws.Range("B7:D7").Copy()
$ws.Range("B8:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row content: "Base flow period starts" criterion with two equal
# date values (2021-09-15).
$ws.Range("A8").Value = "Base flow period starts"
$ws.Range("C8").Value = 44454
$ws.Range("D8").Value = 44454

# "Season start"/"Season end" labels (now rows 6 & 7) drop their bold
# styling.
$ws.Range("A6").Font.Bold = $false
$ws.Range("A7").Font.Bold = $false

# Update the active selection to match the new layout.
$ws.Range("G15").Select()
